$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that currently sits right
#    after the H1 title at the top of the document.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Replace the old "Prompt: ..." image-prompt text (now the last
#    paragraph in the document) with the new meta-description text,
#    keeping the existing (italic) run formatting intact.
$oldPrompt = 'Prompt: Create a cartoon-style image to feature in a blog post about the online slot game "Dragon Riches". The image should prominently feature a happy Maya warrior with glasses. The warrior should be holding a golden dragon and surrounded by Chinese lucky charms, such as coins, paper lanterns, and yuanbao ingots. The background should include elements of both Mayan and Chinese culture, such as temples and dragons. The image should give off a cheerful, lucky, and adventurous vibe to entice readers to give the game a try.'
$newDescription = 'Discover Dragon Riches, an oriental-themed slot game with innovative features and free spin mode. Play it now for free and experience the best of Chinese lucky charms and dragons.'
$d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newDescription, 2) | Out-Null

# 3. Insert a new bold paragraph containing the page title right before
#    that paragraph (inserted right after the paragraph before it, so the
#    "Prompt" paragraph's own run structure is left untouched).
$count = $d.Paragraphs.Count
$prevPara = $d.Paragraphs.Item($count - 1)
$prevPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($count)
$promptPara = $d.Paragraphs.Item($count + 1)
$newPara.Style = $promptPara.Style

$titleText = "Play Dragon Riches Free Slot Game | Oriental-Themed and Innovative Features"
$startPos = $newPara.Range.Start
$newPara.Range.Text = $titleText
$titleRange = $d.Range($startPos, $startPos + $titleText.Length)
$titleRange.Font.Bold = -1
